$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a 0/1-looking value as TEXT (shared string), matching the
# existing crit_* columns which store "0"/"1" as strings rather than numbers.
# Building it via a throw-away formula and flattening with copy/paste-values
# avoids any NumberFormat/style churn (keeps default style 0, no quotePrefix).
# ---------------------------------------------------------------------------
function Set-TextFlag($addr, $flag) {
    $ws.Range($addr).Formula = "=""$flag"""
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
}

# ---------------------------------------------------------------------------
# 1. Header row: a new "crit_and_test" column is being inserted logically at
#    F, pushing the old F:I headers (and their cell comments) one slot right
#    to G:J. Row 9/10 summary formulas stay bound to their original letters
#    (per the target state), so this is done as explicit header/comment
#    moves rather than a real column insert.
# ---------------------------------------------------------------------------

# Comments: capture all the existing text up front, then move right-to-left.
$txtF1 = $ws.Range("F1").Comment.Text()
$txtG1 = $ws.Range("G1").Comment.Text()
$txtH1 = $ws.Range("H1").Comment.Text()
$txtI1 = $ws.Range("I1").Comment.Text()

$ws.Range("I1").Comment.Delete()
$ws.Range("J1").AddComment($txtI1)

$ws.Range("H1").Comment.Delete()
$ws.Range("I1").AddComment($txtH1)

$ws.Range("G1").Comment.Delete()
$ws.Range("H1").AddComment($txtG1)

$ws.Range("F1").Comment.Delete()
$ws.Range("G1").AddComment($txtF1)

# New comment for the inserted F1 header.
$ws.Range("F1").AddComment("test criterio AND")

# Header values: capture up front (Value2, not Value - Value is unreliable
# for reads in this runtime), then move right-to-left, same as the comments.
$hF1 = $ws.Range("F1").Value2
$hG1 = $ws.Range("G1").Value2
$hH1 = $ws.Range("H1").Value2
$hI1 = $ws.Range("I1").Value2

$ws.Range("J1").Value2 = $hI1
$ws.Range("I1").Value2 = $hH1
$ws.Range("H1").Value2 = $hG1
$ws.Range("G1").Value2 = $hF1
$ws.Range("F1").Value2 = "crit_and_test"

# ---------------------------------------------------------------------------
# 2. Data rows: the "preg_test_1" indicator column (numeric 0/1) moves from
#    F to G; the new "crit_and_test" text column is populated at F.
# ---------------------------------------------------------------------------
$f3 = $ws.Range("F3").Value2   # 0.0
$f4 = $ws.Range("F4").Value2   # 1.0

$ws.Range("G3").Value2 = $f3
$ws.Range("G4").Value2 = $f4
$ws.Range("G5").Value2 = 1.0   # new row

Set-TextFlag "F3" "0"
Set-TextFlag "F4" "1"
Set-TextFlag "F5" "0"

# Newly-populated crit_edad / crit_edad_avg cells for rows 4 and 5.
Set-TextFlag "D4" "1"
Set-TextFlag "E4" "0"
Set-TextFlag "D5" "0"
Set-TextFlag "E5" "1"

# ---------------------------------------------------------------------------
# 3. New summary formulas for the appended "pregunta_test_3" column (J).
# ---------------------------------------------------------------------------
$ws.Range("J9").Formula = "=COUNTIF(J1:J7,1)"
$ws.Range("J10").Formula = "=COUNTIF(J1:J7,0)"

$excel.CutCopyMode = 0
